# CIERRE 21 MAY 22
# Roll the payroll sheet forward one week: update the week-label string,
# the hours/extra amount, and the descuento amount. The TODAY()-driven
# date cells and the SUM() totals recalc on their own since they are
# formulas, not literals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week label (B9 feeds H9/B27/H27/B43 through formulas that reference it)
$ws.Range("B9").Value = "SEMANA   20  DEL    16      Al   22   DE   MAYO          2022"

# Horas/extra amount for the week changed 1120 -> 1680 (K24 totals via SUM)
$ws.Range("K21").Value = 1680

# Descuento amount changed 1250 -> 0 (E41 totals via SUM)
$ws.Range("E40").Value = 0

# Move the view/selection from H41 to E41, scrolled down to row 25.
$null = $ws.Range("A25").Select()
$null = $ws.Range("E41").Select()
